$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New header cell "Save" in H1, matching the style of the existing header row
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# New data values in H2:H3
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
